$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 52
$ws.Range("I2").Value = 152
$ws.Range("J2").Value = 573
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 151
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 100
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 65
$ws.Range("T2").Value = 93
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 938
$ws.Range("X2").Value = 870
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 11
